# Actualización automática 2025-07-24 14:55:08
#
# Updates the monthly sales figures for salesperson "CARRION CARRION LESLY
# ANABE" (row 5 on the first two sheets) and propagates the new totals /
# compliance ratios that depend on them on the "CUMPLIMIENTO MENSUAL" sheet.

$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual     = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento     = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- VENTAS POR GRUPO ---------------------------------------------------
# M5: updated monthly sales amount for CARRION CARRION LESLY ANABE.
$wsVentasPorGrupo.Range("M5").Value = 7043.75

# --- VENTA MENSUAL -------------------------------------------------------
# F5: updated monthly sales amount for CARRION CARRION LESLY ANABE.
$wsVentaMensual.Range("F5").Value = 7445.56

# F22: column total (sum of F4:F21) recalculated after the F5 change.
$wsVentaMensual.Range("F22").Value = 42599.07

# --- CUMPLIMIENTO MENSUAL ------------------------------------------------
# Row 16 ("PORCELANATO"): achieved amount (D), remaining gap (E) and
# achievement ratio (F = D/C) recalculated from the new sales figures.
$wsCumplimiento.Range("D16").Value = 33386.02
$wsCumplimiento.Range("E16").Value = 10880.22
$wsCumplimiento.Range("F16").Value = 0.7542095285255761

# Row 19 ("TOTAL"): grand totals recalculated to reflect the row 16 change.
$wsCumplimiento.Range("D19").Value = 42599.06999999999
$wsCumplimiento.Range("E19").Value = 22778.92762291768
$wsCumplimiento.Range("F19").Value = 0.6515811366034751
